$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-19 Friday" "2025-09-20 Saturday"

Replace-Text "869÷4=" "293÷7="
Replace-Text "806÷4=" "331÷7="
Replace-Text "835÷7=" "532÷4="
Replace-Text "147÷4=" "982÷9="
Replace-Text "937÷8=" "980÷7="
Replace-Text "126÷4=" "657÷2="
Replace-Text "119÷5=" "649÷9="
Replace-Text "104÷5=" "226÷7="
Replace-Text "968÷7=" "761÷6="
Replace-Text "956÷9=" "889÷5="
Replace-Text "192÷5=" "409÷5="
Replace-Text "171÷9=" "217÷6="
Replace-Text "522÷3=" "258÷9="
Replace-Text "346÷4=" "672÷9="
Replace-Text "685÷2=" "244÷6="
Replace-Text "850÷2=" "554÷5="
Replace-Text "879÷4=" "321÷3="
Replace-Text "522÷7=" "816÷5="
Replace-Text "985÷6=" "443÷3="
Replace-Text "621÷2=" "260÷4="
Replace-Text "179÷9=" "203÷4="
Replace-Text "261÷6=" "286÷4="
Replace-Text "656÷6=" "718÷7="
Replace-Text "713÷2=" "443÷6="
Replace-Text "932÷9=" "146÷5="
